# Refresh the scraped cryptos list: update Price (D) / Volume(1h) (E) for
# most rows, plus the dogwifhat/ThetaToken rank swap in rows 43-44 (B/C/D/E).
# Values in D that look like plain numbers are written with a leading
# apostrophe so Excel stores them as text (matching the original
# inlineStr cells, e.g. "593.23") instead of silently coercing them to a
# Double and mangling values like "70.756.37" or trailing zeros like "9.90".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.756.37"
$ws.Range("E2").Value = "  +5.83%  "
$ws.Range("D3").Value = "3.636.89"
$ws.Range("E3").Value = "  +5.69%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'593.23"
$ws.Range("E5").Value = "  +1.99%  "
$ws.Range("D6").Value = "'194.90"
$ws.Range("E6").Value = "  +2.95%  "
$ws.Range("D7").Value = "'0.646"
$ws.Range("E7").Value = "  +2.71%  "
$ws.Range("D8").Value = "3.630.92"
$ws.Range("E8").Value = "  +5.72%  "
$ws.Range("D9").Value = "'0.999"
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("E10").Value = "  +8.57%  "
$ws.Range("E11").Value = "  +5.58%  "
$ws.Range("D12").Value = "'58.07"
$ws.Range("E12").Value = "  +1.33%  "
$ws.Range("D13").Value = "'0.0000306"
$ws.Range("E13").Value = "  +10.72%  "
$ws.Range("D14").Value = "'9.99"
$ws.Range("E14").Value = "  +5.67%  "
$ws.Range("D15").Value = "4.222.18"
$ws.Range("E15").Value = "  +5.70%  "
$ws.Range("D16").Value = "'20.48"
$ws.Range("E16").Value = "  +8.78%  "
$ws.Range("D17").Value = "3.637.79"
$ws.Range("E17").Value = "  +5.56%  "
$ws.Range("D18").Value = "70.767.53"
$ws.Range("E18").Value = "  +5.87%  "
$ws.Range("E19").Value = "  +5.86%  "
$ws.Range("E20").Value = "  +2.77%  "
$ws.Range("D21").Value = "'1.07"
$ws.Range("E21").Value = "  +3.84%  "
$ws.Range("D22").Value = "'488.60"
$ws.Range("E22").Value = "  +2.24%  "
$ws.Range("D23").Value = "'19.34"
$ws.Range("E24").Value = "  -2.31%  "
$ws.Range("E25").Value = "  +3.08%  "
$ws.Range("D26").Value = "'91.43"
$ws.Range("E26").Value = "  +2.42%  "
$ws.Range("E27").Value = "  +6.49%  "
$ws.Range("E28").Value = "  +4.25%  "
$ws.Range("D29").Value = "'9.58"
$ws.Range("E29").Value = "  +6.31%  "
$ws.Range("D30").Value = "'7.93"
$ws.Range("E30").Value = "  +6.48%  "
$ws.Range("D31").Value = "'32.82"
$ws.Range("E31").Value = "  +5.37%  "
$ws.Range("E32").Value = "  +9.92%  "
$ws.Range("D33").Value = "'12.29"
$ws.Range("E33").Value = "  +4.45%  "
$ws.Range("D34").Value = "'616.78"
$ws.Range("E34").Value = "  +2.56%  "
$ws.Range("D35").Value = "'66.39"
$ws.Range("E35").Value = "  +2.97%  "
$ws.Range("D36").Value = "'40.21"
$ws.Range("E36").Value = "  +8.12%  "
$ws.Range("D37").Value = "0.0₃0832"
$ws.Range("E37").Value = "  +10.90%  "
$ws.Range("E38").Value = "  +5.75%  "
$ws.Range("E39").Value = "  +1.45%  "
$ws.Range("D40").Value = "'0.998"
$ws.Range("E40").Value = "  -0.12%  "
$ws.Range("E41").Value = "  +2.51%  "
$ws.Range("D42").Value = "3.326.17"
$ws.Range("E42").Value = "  +4.04%  "
$ws.Range("B43").Value = "ThetaToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D43").Value = "'3.17"
$ws.Range("E43").Value = "  +8.77%  "
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "'3.17"
$ws.Range("E44").Value = "  +16.00%  "
$ws.Range("E45").Value = "  +8.16%  "
$ws.Range("D46").Value = "'0.0458"
$ws.Range("E46").Value = "  +6.36%  "
$ws.Range("D47").Value = "'9.65"
$ws.Range("E47").Value = "  +11.39%  "
$ws.Range("E48").Value = "  +2.71%  "
$ws.Range("D49").Value = "'0.139"
$ws.Range("E49").Value = "  +3.48%  "
$ws.Range("E50").Value = "  +1.62%  "
$ws.Range("E51").Value = "  +0.13%  "
